$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New loading_percent values for rows 2-25 (columns B,C,E,F,G,H,I,J,L)
$data = @{
    2 = @(18.83950662777261, 9.699123079810191, 14.43083828962585, 38.78817939157862, 37.80507229859841, 16.44182488177383, 25.3368332873355, 8.014051268546291, 12.8381308817174)
    3 = @(18.27168214458624, 9.314517981662007, 14.44789008494764, 38.83513334338183, 37.86096500995178, 16.51038522740477, 25.48689780785287, 8.03311775498371, 12.80276705270379)
    4 = @(17.91632557647129, 9.068318746583692, 14.45957340809967, 38.8752462922676, 37.91224727598833, 16.55655223207637, 25.58590765777046, 8.045443115289366, 12.7826101563789)
    5 = @(17.77004374688135, 8.965550215830758, 14.46464021746733, 38.89442112168817, 37.9373823618694, 16.57638575029594, 25.62797706989208, 8.050621851208033, 12.77479249283928)
    6 = @(17.74567152563751, 8.948340958002149, 14.46550004271552, 38.89777563676509, 37.94181103129729, 16.5797405983105, 25.63506652027088, 8.051491217228723, 12.77351845982803)
    7 = @(17.91435842824095, 9.066942531695938, 14.45964050201226, 38.87549345103287, 37.91256914107542, 16.55681558826062, 25.58646805492252, 8.045512324989041, 12.78250311309023)
    8 = @(18.64523827232014, 9.568645696234393, 14.43646627549446, 38.80202250125284, 37.82080717839803, 16.46461751797581, 25.38714651358355, 8.020497353421764, 12.82561744233203)
    9 = @(20.01655544115571, 10.46949251705491, 14.40062237305043, 38.74779868166156, 37.77655931372693, 16.31627216988151, 25.05102654293668, 7.976325986459259, 12.92228456170398)
    10 = @(20.97580124398428, 11.07711177478106, 14.38010381159231, 38.76308839701899, 37.82799855794583, 16.22727920835938, 24.83779161832284, 7.946816182886926, 13.0003689945954)
    11 = @(21.39980876027322, 11.34116513756179, 14.37202427645658, 38.78205585810753, 37.86980017860772, 16.19118044681542, 24.74818069610386, 7.93402325517472, 13.03735031648044)
    12 = @(21.55845483238572, 11.43933880733508, 14.36914447402175, 38.7909661035754, 37.88828251761252, 16.17814489077865, 24.71531669199151, 7.929269131832976, 13.05155737182338)
    13 = @(21.52437485563659, 11.41827677212765, 14.36975670736636, 38.78897029704609, 37.88418396418582, 16.18092406278949, 24.7223468559103, 7.930289009833673, 13.0484887117658)
    14 = @(21.4128998207766, 11.34927858322717, 14.3717837548278, 38.78275428838749, 37.87126753962335, 16.1900952708991, 24.74545547598927, 7.933630324012266, 13.03851511114527)
    15 = @(21.34436469184922, 11.30677739154523, 14.37304876824004, 38.77917177865104, 37.86370147547723, 16.19579561087152, 24.75974968856881, 7.935688718200347, 13.03243222276537)
    16 = @(20.94783043377668, 11.05960299417577, 14.38065701607261, 38.7620906137881, 37.82563770250306, 16.22972689950048, 24.84379721227906, 7.94766489081448, 12.99798097413437)
    17 = @(20.70130046405679, 10.90477372166385, 14.3856452700765, 38.75468927350362, 37.80700625132049, 16.25166834030726, 24.89725547292985, 7.955173208167174, 12.97721539735589)
    18 = @(20.5583470924232, 10.81455903941802, 14.38863249718666, 38.75156294466295, 37.79802186814121, 16.26470095320885, 24.9286984765677, 7.959551237034238, 12.9654094717428)
    19 = @(20.50975150462707, 10.78381585271293, 14.38966422692511, 38.75069858674074, 37.79527708730294, 16.26918429948992, 24.93946373403946, 7.961043786810171, 12.96143607814768)
    20 = @(20.7276646823406, 10.92137603575007, 14.38510204208537, 38.75536012496088, 37.80881028852183, 16.24928991939861, 24.89149274701571, 7.954367786400454, 12.97941170754887)
    21 = @(21.4456957878715, 11.3695946469766, 14.37118348892337, 38.78453319806962, 37.87498937427758, 16.1873842214042, 24.73863882663115, 7.932646452724658, 13.04143914586207)
    22 = @(21.90374609115453, 11.65192244601928, 14.36313427585622, 38.81366863171551, 37.93370466302361, 16.15062393077305, 24.64497834002962, 7.918976307130176, 13.08315769053606)
    23 = @(21.66034533703717, 11.50222137849804, 14.36733467534928, 38.79719751848327, 37.90095124933806, 16.16990393518688, 24.69439356596698, 7.926224352361726, 13.0607861060892)
    24 = @(20.71574920857372, 10.91387386813594, 14.38534726346848, 38.75505331698324, 37.80798930420543, 16.25036390115416, 24.89409586701311, 7.954731726290378, 12.9784183425353)
    25 = @(19.6533558746503, 10.2350858380232, 14.40929504839355, 38.75280558927003, 37.77385937728253, 16.35290704576883, 25.13606046379489, 7.987756279948553, 12.89486848136319)
}

$cols = @(2, 3, 5, 6, 7, 8, 9, 10, 12)

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Cells.Item($row, $cols[$i]).Value = $values[$i]
    }
}
